# Documentation/Specifications.docx update:
# "Project version" attribute type documentation was "any Integer*" with a
# footnote ("* = May change this to a String soon"); the type actually
# changed to String a while ago, so the spec doc is corrected to say
# "String" outright and the now-obsolete footnote text is removed.

$d = $word.ActiveDocument

# 1) "version" attribute's type cell: "any Integer*" -> "String"
$d.Content.Find.Execute("any Integer*", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "String", 2) | Out-Null

# 2) Remove the obsolete footnote paragraph text, leaving an empty paragraph
#    (the paragraph mark / spacing stays, only the run with the note is gone).
$r = $d.Content
$r.Find.Execute("* = May change this to a String soon", $false, $false, $false, `
                 $false, $false, $true, 1, $false, "", 2) | Out-Null

# Word always keeps a "_GoBack" bookmark marking the last edited spot; move it
# (by re-adding it, which replaces any existing bookmark of the same name) to
# the location of this last edit, matching what real Word does when you type
# there.
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
